$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, bordered, centered style) from H1 to the
# two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 3

$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 8

$ws.Range("I8").Value = 8
$ws.Range("J8").Value = 8

$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 6

$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 4
